# Auto-generated edit script applying the cryptos.xlsx data refresh diff
# (GitHub Actions "Updated cryptos list" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.367.07"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "2.005.60"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.00%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.852"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.68%  "
$ws.Range("D15").Value = "2.300.63"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "2.007.69"
$ws.Range("E17").Value = "  +3.76%  "
$ws.Range("D18").Value = "37.271.80"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").Value = "0.0₃0868"
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.37%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.143"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.39%  "
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.33%  "
$ws.Range("E35").Value = "  +7.50%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.76%  "
$ws.Range("D46").Value = "1.380.37"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("E47").Value = "  +3.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +17.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.03%  "
